$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'26.793.40"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +1.69%  "
$ws.Range("D3").Value = "'1.720.95"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.45%  "
$ws.Range("D4").Value = "'0.9994"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.32%  "
$ws.Range("D5").Value = "'239.77"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.29%  "
$ws.Range("D6").Value = "'0.9997"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.25%  "
$ws.Range("D7").Value = "'0.4749"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -2.30%  "
$ws.Range("D8").Value = "'0.2553"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.13%  "
$ws.Range("E9").Value = "  -0.95%  "
$ws.Range("D10").Value = "'1.719.59"
$ws.Range("D10").Style = "Normal"
$ws.Range("D11").Value = "'15.80"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +2.14%  "
$ws.Range("D12").Value = "'0.06892"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.81%  "
$ws.Range("D13").Value = "'0.5942"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.41%  "
$ws.Range("D14").Value = "'4.391"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.68%  "
$ws.Range("D15").Value = "'76.19"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.26%  "
$ws.Range("E16").Value = "  +0.32%  "
$ws.Range("D17").Value = "'26.706.49"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.77%  "
$ws.Range("D18").Value = "'0.9994"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.31%  "
$ws.Range("D19").Value = "'0.000006987"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.63%  "
$ws.Range("D20").Value = "'11.25"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.28%  "
$ws.Range("D21").Value = "'1.940.12"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.31%  "
$ws.Range("D22").Value = "'4.363"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.94%  "
$ws.Range("D23").Value = "'8.310"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.29%  "
$ws.Range("D24").Value = "'5.052"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.26%  "
$ws.Range("D25").Value = "'140.75"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +3.46%  "
$ws.Range("D26").Value = "'15.10"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.28%  "
$ws.Range("E27").Value = "  +3.22%  "
$ws.Range("B28").Value = "BitcoinCash"
$ws.Range("C28").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D28").Value = "'105.92"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.63%  "
$ws.Range("B29").Value = "Toncoin"
$ws.Range("C29").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D29").Value = "'1.371"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.69%  "
$ws.Range("E30").Value = "  +1.58%  "
$ws.Range("D31").Value = "'0.07875"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.77%  "
$ws.Range("D32").Value = "'3.625"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.61%  "
$ws.Range("D33").Value = "'0.04607"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +4.25%  "
$ws.Range("D34").Value = "'2.591"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.21%  "
$ws.Range("D35").Value = "'0.9913"
$ws.Range("D35").Style = "Normal"
$ws.Range("D36").Value = "'0.6081"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.33%  "
$ws.Range("D37").Value = "'0.9166"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.96%  "
$ws.Range("D38").Value = "'2.501"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +5.56%  "
$ws.Range("E39").Value = "  -0.61%  "
$ws.Range("D40").Value = "'0.9989"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.28%  "
$ws.Range("E41").Value = "  +5.07%  "
$ws.Range("E42").Value = "  +0.29%  "
$ws.Range("D43").Value = "'99.58"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.03%  "
$ws.Range("D44").Value = "'0.3774"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.64%  "
$ws.Range("D45").Value = "'6.696"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.85%  "
$ws.Range("D46").Value = "'0.1141"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.61%  "
$ws.Range("D47").Value = "'0.05339"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.02%  "
$ws.Range("D48").Value = "'7.694"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.11%  "
$ws.Range("D49").Value = "'29.66"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -3.04%  "
$ws.Range("D50").Value = "'1.231"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.80%  "
$ws.Range("D51").Value = "'1.001"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.25%  "
